# "balance the weapon and skills"
# Add a new Missile row (Id=38: stone / 飞石 / yellowsplash) to the Missile sheet's
# data table, growing the table/autofilter/dimension from row 40 to row 41,
# widening the effect-stat columns (E:H) and updating the view/selection to
# point at the freshly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write the new data row (row 41) that the table will grow to include.
#    Id | TypeName | Name | EffName | Speed | Image | FrameCount | FrameTime
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = 38
$ws.Range("B41").Value = "stone"
$ws.Range("C41").Value = "飞石"
$ws.Range("D41").Value = "yellowsplash"
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 370
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 4

# ---------------------------------------------------------------------------
# 2. Grow the worksheet Table ("表1") / AutoFilter so it spans the new row.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:H41"))

# ---------------------------------------------------------------------------
# 3. Give the new stat columns (E:H) an explicit width, matching the other
#    formatted columns in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("E1:H1").ColumnWidth = 5.160714285714286

# ---------------------------------------------------------------------------
# 4. Update the window/view so the newly-added row is visible and selected,
#    the same way Excel leaves the view after you scroll down and fill in a
#    new table row.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$excel.Goto($ws.Range("D41"), $true)

Write-Host "Missile.xlsx: added row 38 (stone/yellowsplash) and resized table to A3:H41"
